# Balance Sheet Report template - "final version"
#
# Adds ASSETS / LIABILITIES / EQUITY sections (each: a blue "Heading 1"
# banner, a "{{item.AccountName}}" / "{{item.Balance}}" placeholder row, and
# a bold "TOTAL ..." row with a currency <<sum>> placeholder finished off
# with the built-in "Total" cell style), a grand "TOTAL LIABILITIES &
# EQUITY" row with a real formula, and three workbook-scoped named ranges
# (assets / liabilities / equity) pointing at the placeholder+total rows of
# each section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$CURRENCY_FMT = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

$xlPasteFormats = -4122
$xlDouble = -4119
$xlRight = -4152
$xlLeft = -4131

# ---------------------------------------------------------------------------
# Defined names (workbook scope) - one per section, covering the
# placeholder row + the TOTAL row underneath it.
# ---------------------------------------------------------------------------
$wb.Names.Add("assets", "=Sheet1!`$A`$8:`$G`$9")
$wb.Names.Add("liabilities", "=Sheet1!`$A`$14:`$G`$15")
$wb.Names.Add("equity", "=Sheet1!`$A`$20:`$G`$21")

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Section banner row (e.g. "ASSETS"/"LIABILITIES"/"EQUITY"), matching the
# pre-existing "ASSETS" row's built-in "Heading 1" cell style.
function Add-SectionHeading([int]$row, [string]$text) {
    $rng = $ws.Range("A$($row):G$($row)")
    $rng.Merge() | Out-Null
    $rng.Value = $text
    $rng.Style = "Heading 1"
    $ws.Rows.Item($row).RowHeight = 20.25
}

# Item placeholder row: "{{item.AccountName}}" (A:F, plain/general) and
# "{{item.Balance}}" (G, currency, right-aligned, plain font/border).
function Add-PlaceholderRow([int]$row) {
    $labelRng = $ws.Range("A$($row):F$($row)")
    $labelRng.Merge() | Out-Null
    $labelRng.Value = "{{item.AccountName}}"

    $amountCell = $ws.Range("G$($row)")
    $amountCell.Value = "{{item.Balance}}"
    $amountCell.NumberFormat = $CURRENCY_FMT
    $amountCell.HorizontalAlignment = $xlRight
}

# TOTAL row for a section, e.g. "TOTAL ASSETS" (A:F) + "<<sum>>" (G),
# finished off with the built-in "Total" cell style (double-rule under a
# thin top rule). $bigFont selects the larger, section-header-matching
# font used for the very first "TOTAL ASSETS" row in the template.
function Add-TotalRow([int]$row, [string]$label, [bool]$bigFont) {
    $labelRng = $ws.Range("A$($row):F$($row)")
    $labelRng.Merge() | Out-Null
    $labelRng.Value = $label
    $labelRng.Style = "Total"

    $amountCell = $ws.Range("G$($row)")
    $amountCell.Value = "<<sum>>"
    $amountCell.Style = "Total"
    $amountCell.NumberFormat = $CURRENCY_FMT
    $amountCell.HorizontalAlignment = $xlRight

    if ($bigFont) {
        $labelRng.Font.Size = 14
        $labelRng.Font.ThemeColor = 3
    }
}

# ---------------------------------------------------------------------------
# G5 "TOTALS" column header -- same font family as the section headings
# (bold 14pt, theme accent) but right-aligned, no border.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$totalsHeader = $ws.Range("G5")
$totalsHeader.PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$totalsHeader.Value = "TOTALS"
$totalsHeader.HorizontalAlignment = $xlRight
$ws.Rows.Item(5).RowHeight = 18.75

# ---------------------------------------------------------------------------
# ASSETS section (rows 8-9)
# ---------------------------------------------------------------------------
Add-PlaceholderRow 8
Add-TotalRow 9 "TOTAL ASSETS" $true
$ws.Rows.Item(9).RowHeight = 19.5
$ws.Rows.Item(10).RowHeight = 15.75

# ---------------------------------------------------------------------------
# LIABILITIES section (rows 12, 14-15)
# ---------------------------------------------------------------------------
Add-SectionHeading 12 "LIABILITIES"
$ws.Rows.Item(13).RowHeight = 15.75
Add-PlaceholderRow 14
Add-TotalRow 15 "TOTAL LIABILITIES" $false
$ws.Rows.Item(15).RowHeight = 15.75
$ws.Rows.Item(16).RowHeight = 15.75

# ---------------------------------------------------------------------------
# EQUITY section (rows 18, 20-21)
# ---------------------------------------------------------------------------
Add-SectionHeading 18 "EQUITY"
$ws.Rows.Item(19).RowHeight = 15.75
Add-PlaceholderRow 20
Add-TotalRow 21 "TOTAL EQUITY" $false
$ws.Rows.Item(21).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Grand total row (rows 22 spacer, 23 total, 24 spacer)
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).RowHeight = 16.5

$grandLabel = $ws.Range("A23:F23")
$grandLabel.Merge() | Out-Null
$grandLabel.Value = "TOTAL LIABILITIES & EQUITY"
$grandLabel.Style = "Total"
$grandTop = $grandLabel.Borders.Item(8)
$grandTop.LineStyle = $xlDouble
$grandLabel.Font.Size = 14
$grandLabel.Font.ThemeColor = 3
$grandLabel.HorizontalAlignment = $xlLeft

$grandAmount = $ws.Range("G23")
$grandAmount.Style = "Total"
$grandAmountTop = $grandAmount.Borders.Item(8)
$grandAmountTop.LineStyle = $xlDouble
$grandAmount.NumberFormat = $CURRENCY_FMT
$grandAmount.HorizontalAlignment = $xlRight
$grandAmount.Formula = "=G`$15+G`$21"

$ws.Rows.Item(23).RowHeight = 20.25
$ws.Rows.Item(24).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Selection, matching the template's final saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("A15:F15").Select() | Out-Null
